# Remove a handful of Event Points that were determined to be "on the day side".
# These correspond (in the original row numbering) to rows 10, 24, 30 and 34 of
# Sheet1 - each one a full record (Narrowed Point date/time, USED SAT, General
# Point date/time, Kp, F10.7, Interval Start/End, Interval Kp/F10.7, Notes).
#
# Deleting them from top to bottom causes each subsequent row to shift up, so
# the row indices used below account for that shift as we go.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Rows.Item(10).Delete()
$ws.Rows.Item(23).Delete()
$ws.Rows.Item(28).Delete()
$ws.Rows.Item(31).Delete()
